$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster, Ligand symbol, Receptor symbol, Target cluster, and the 16
# numeric TPM-derived metric columns for every (sending, target) cluster pair.
$newData = @(
    ,("ECs", "Angptl1", "Tek", "ECs", 1, 0.3333333333333333, 0.008964666666666668, 0.026894, 0.0001537549234563009, 0.0001537549234563009, 3, 1, 43.28265566666667, 129.847967, 0.667219228070094, 0.667219228070094, 0.3880145804997778, 3.492131224498, 0.0001025882413404895, 0.0001025882413404895)
    ,("ECs", "Angptl1", "Tek", "FAPs", 1, 0.3333333333333333, 0.008964666666666668, 0.026894, 0.0001537549234563009, 0.0001537549234563009, 3, 1, 18.63243533333333, 55.897306, 0.2872263480299067, 0.2872263480299067, 0.1670335719515556, 1.503302147564, 0.00004416246515597115, 0.00004416246515597114)
    ,("ECs", "Angptl1", "Tek", "MuSCs", 1, 0.3333333333333333, 0.008964666666666668, 0.026894, 0.0001537549234563009, 0.0001537549234563009, 3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 0.02649171357155556, 0.238425422144, 0.000007004216959840299, 0.000007004216959840298)
    ,("FAPs", "Angptl1", "Tek", "ECs", 3, 1, 55.82849233333334, 167.485477, 0.9575264630094835, 0.9575264630094834, 3, 1, 43.28265566666667, 129.847967, 0.667219228070094, 0.667219228070094, 2416.405410052807, 21747.64869047526, 0.638880067505875, 0.638880067505875)
    ,("FAPs", "Angptl1", "Tek", "FAPs", 3, 1, 55.82849233333334, 167.485477, 0.9575264630094835, 0.9575264630094834, 3, 1, 18.63243533333333, 55.897306, 0.2872263480299067, 0.2872263480299067, 1040.220773158329, 9361.986958424963, 0.2750268291122074, 0.2750268291122074)
    ,("FAPs", "Angptl1", "Tek", "MuSCs", 3, 1, 55.82849233333334, 167.485477, 0.9575264630094835, 0.9575264630094834, 3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 164.9801920160391, 1484.821728144352, 0.04361956639140114, 0.04361956639140114)
    ,("MuSCs", "Angptl1", "Tek", "ECs", 3, 1, 2.467451, 7.402353, 0.0423197820670603, 0.04231978206706029, 3, 1, 43.28265566666667, 129.847967, 0.667219228070094, 0.667219228070094, 106.7978320073723, 961.180488066351, 0.02823657232287858, 0.02823657232287857)
    ,("MuSCs", "Angptl1", "Tek", "FAPs", 3, 1, 2.467451, 7.402353, 0.0423197820670603, 0.04231978206706029, 3, 1, 18.63243533333333, 55.897306, 0.2872263480299067, 0.2872263480299067, 45.97462119566867, 413.771590761018, 0.01215535645254326, 0.01215535645254326)
    ,("MuSCs", "Angptl1", "Tek", "MuSCs", 3, 1, 2.467451, 7.402353, 0.0423197820670603, 0.04231978206706029, 3, 1, 2.955125333333334, 8.865376000000001, 0.04555442389999943, 0.04555442389999944, 7.291626958858667, 65.624642629728, 0.001927853291638459, 0.001927853291638459)
)

for ($r = 0; $r -lt $newData.Length; $r++) {
    $rowValues = $newData[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}
